# Updates the "레벨" log sheet: refreshes the B/C values of the four
# existing rows and appends seven new rows (rows 5-11), matching the
# "Add files via upload" commit.
#
# Column A holds purely-numeric, 18-digit snowflake-style ID strings.
# Excel's COM layer auto-converts plain numeric-looking strings assigned
# to .Value into Number values, which silently loses precision for
# integers this large (they'd round-trip as 7.68...E+17 instead of the
# exact digit string). Setting NumberFormat to "@" (Text) before the
# assignment - the moral equivalent of Excel's manual quote-prefix entry
# - keeps these values as exact text, matching the original inlineStr
# cells in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "523017072796499968"; B = -1;  C = "2021/01/11 22:22:03" },
    @{ A = "418022156438601738"; B = 64;  C = "2021/01/11 21:18:01" },
    @{ A = "720392147928350812"; B = 114; C = "2021/01/11 22:01:55" },
    @{ A = "450168301369163786"; B = 5;   C = "2021/01/11 21:05:51" },
    @{ A = "768086907317649430"; B = 42;  C = "2021/01/11 22:22:25" },
    @{ A = "441003067668955166"; B = 259; C = "2021/01/11 21:05:51" },
    @{ A = "508195540245151762"; B = 13;  C = "2021/01/11 21:54:34" },
    @{ A = "523017072796499968"; B = -2;  C = "2021/01/11 21:43:02" },
    @{ A = "523017072796499968"; B = 10;  C = "2021/01/11 21:43:04" },
    @{ A = "523017072796499968"; B = -3;  C = "2021/01/11 21:43:06" },
    @{ A = "523017072796499968"; B = 5;   C = $null }
)

$r = 1
foreach ($row in $rows) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    if ($row.C -ne $null) {
        $ws.Range("C$r").Value = $row.C
    }
    $r = $r + 1
}
